$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: 展览 (Exhibition)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F6").Value  = 4327
$ws1.Range("F14").Value = 642
$ws1.Range("F15").Value = 1531
$ws1.Range("F16").Value = 1414
$ws1.Range("F18").Value = 29
$ws1.Range("F19").Value = 557
$ws1.Range("F20").Value = 4034
$ws1.Range("G20").Value = 80
$ws1.Range("F21").Value = 4034
$ws1.Range("G21").Value = 80
$ws1.Range("F22").Value = 661
$ws1.Range("F24").Value = 766
$ws1.Range("F26").Value = 2223
$ws1.Range("F30").Value = 34
$ws1.Range("F31").Value = 1177
$ws1.Range("F33").Value = 51
$ws1.Range("F34").Value = 1070
$ws1.Range("F35").Value = 1080

# ---------------------------------------------------------------------------
# Sheet 2: 演出 (Performance)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F6").Value  = 97
$ws2.Range("F18").Value = 266

# ---------------------------------------------------------------------------
# Sheet 3: 本地生活 (Local life)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F4").Value = 551
$ws3.Range("F5").Value = 82
$ws3.Range("F6").Value = 176

# ---------------------------------------------------------------------------
# Sheet 4: 全部类型 (All types)
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

# Row 7 now holds what used to be row 8's event (银魂同人only同好会1.0)
$ws4.Range("B7").Value = "'2024-10-13"
$ws4.Range("C7").Value = "北京·银魂同人only同好会1.0"
$ws4.Range("D7").Value = "丽泽天地购物中心 丽泽天地购物中心"
$ws4.Range("E7").Value = "2024.10.13 10:00-10.13 18:00"
$ws4.Range("F7").Value = 82
$ws4.Range("G7").Value = 98
$ws4.Range("H7").Value = "https://show.bilibili.com/platform/detail.html?id=93073"
$ws4.Range("I7").Value = "//i0.hdslb.com/bfs/openplatform/202409/SqTqdBYb1727607319860.jpeg"

# Row 8 now holds what used to be row 9's event (EVANGELION x PrismLand)
$ws4.Range("B8").Value = "'2024-10-15"
$ws4.Range("C8").Value = "北京·EVANGELION× PrismLand · 新世纪福音战士官方授权主题店"
$ws4.Range("D8").Value = "王府井地铁站F1东口步行120米 北京王府井喜悦购物中心"
$ws4.Range("E8").Value = "2024.10.15 00:00-12.15 23:59"
$ws4.Range("F8").Value = 551
$ws4.Range("G8").Value = 20
$ws4.Range("H8").Value = "https://show.bilibili.com/platform/detail.html?id=93039"
$ws4.Range("I8").Value = "//i0.hdslb.com/bfs/openplatform/202409/n32CfRya1727584778969.jpeg"

# Row 9 is a brand new event (全职高手 x HAPPY ZOO coffee shop)
$ws4.Range("B9").Value = "'2024-10-17"
$ws4.Range("C9").Value = "北京·全职高手×HAPPY ZOO 全职高手十周年咖啡厅"
$ws4.Range("D9").Value = "学清路38号金码大厦B座(六道口地铁站B东北口步行110米) BOM嘻番里"
$ws4.Range("E9").Value = "2024.10.17 00:00-11.17 23:59"
$ws4.Range("F9").Value = 82
$ws4.Range("G9").Value = 10
$ws4.Range("H9").Value = "https://show.bilibili.com/platform/detail.html?id=93324"
$ws4.Range("I9").Value = "//i2.hdslb.com/bfs/openplatform/202410/bVeTwUWu1728699225130.png"

$ws4.Range("F10").Value = 4331
$ws4.Range("F11").Value = 4332
$ws4.Range("F17").Value = 97
$ws4.Range("F25").Value = 1531
$ws4.Range("F27").Value = 1414
$ws4.Range("F29").Value = 557
$ws4.Range("F31").Value = 4034
$ws4.Range("G31").Value = 80
$ws4.Range("F32").Value = 4034
$ws4.Range("G32").Value = 80
$ws4.Range("F33").Value = 661
$ws4.Range("F35").Value = 766
$ws4.Range("F37").Value = 2223
$ws4.Range("F41").Value = 34
$ws4.Range("F42").Value = 1177
$ws4.Range("F44").Value = 266
$ws4.Range("F48").Value = 51
$ws4.Range("F49").Value = 1070
$ws4.Range("F50").Value = 1080
